$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Lake Florida survey rows (235-265): copy number formats/styles from the
# last existing data row (234), then fill in the values cell-by-cell so the
# shared-string table is reused (dedup) exactly as the source workbook does.
$lastRow = 234
$firstNew = 235
$lastNew = 265

$srcFmt = $ws.Range("A" + $lastRow + ":J" + $lastRow)
$dstFmt = $ws.Range("A" + $firstNew + ":J" + $lastNew)
$srcFmt.Copy()
$dstFmt.PasteSpecial(-4122)

$rows = @(
    @(235, 43328.515668078704, 'Double observer distance', 'Florida', 43326.0, 'Aislyn', 13.0, 0.0, 1.1, 'Silt, Sand, Plant cover', '2) 26-50% cover'),
    @(236, 43328.51595217593, 'Double observer distance', 'Florida', 43326.0, 'Aislyn', 13.0, 26.0, 1.9, 'Silt, Sand, Plant cover', '3) 51-75% cover'),
    @(237, 43328.51623849537, 'Double observer distance', 'Florida', 43326.0, 'Aislyn', 14.0, 0.0, 1.2, 'Silt, Sand, Plant cover', '1) 1-25% cover'),
    @(238, 43328.51662375, 'Double observer distance', 'Florida', 43326.0, 'Aislyn', 14.0, 16.0, 1.8, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(239, 43328.51687211805, 'Double observer distance', 'Florida', 43326.0, 'Aislyn', 15.0, 0.0, 1.3, 'Silt, Sand, Gravel, Plant cover', '1) 1-25% cover'),
    @(240, 43328.517160717594, 'Double observer distance', 'Florida', 43326.0, 'Aislyn', 15.0, 14.0, 2.0, 'Silt, Sand, Plant cover', '2) 26-50% cover'),
    @(241, 43328.51741059028, 'Double observer distance', 'Florida', 43326.0, 'Aislyn', 15.0, 20.0, 2.2, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(242, 43328.51898681713, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 11.0, 0.0, 1.0, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(243, 43328.51919377314, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 12.0, 0.0, 1.1, 'Silt, Sand, Plant cover', '3) 51-75% cover'),
    @(244, 43328.51940967592, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 12.0, 8.0, 1.4, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(245, 43328.51965150463, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 13.0, 0.0, 1.2, 'Silt, Sand, Plant cover', '1) 1-25% cover'),
    @(246, 43328.52022525463, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 13.0, 20.0, 1.7, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(247, 43328.520576678246, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 14.0, 0.0, 1.1, 'Silt, Sand, Plant cover', '1) 1-25% cover'),
    @(248, 43328.52086074074, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 14.0, 7.0, 1.5, 'Silt, Sand, Plant cover', '3) 51-75% cover'),
    @(249, 43328.52118605324, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 14.0, 21.0, 2.0, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(250, 43328.52410350695, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 15.0, 0.0, 0.9, 'Silt, Sand, Gravel, Plant cover', '1) 1-25% cover'),
    @(251, 43328.52433186343, 'Double observer no distance', 'Florida', 43327.0, 'Aislyn', 15.0, 17.0, 2.3, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(252, 43328.52475386574, 'Quadrat survey', 'Florida', 43326.0, 'Aislyn', 15.0, 0.0, 1.1, 'Silt, Sand, Gravel, Plant cover', '3) 51-75% cover'),
    @(253, 43328.52498887731, 'Quadrat survey', 'Florida', 43327.0, 'Aislyn', 15.0, 7.0, 1.7, 'Silt, Sand, Plant cover', '1) 1-25% cover'),
    @(254, 43328.52530681713, 'Quadrat survey', 'Florida', 43326.0, 'Aislyn', 15.0, 14.0, 2.0, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(255, 43328.5257092824, 'Quadrat survey', 'Florida', 43326.0, 'Aislyn', 14.0, 0.0, 1.3, 'Silt, Sand, Plant cover', '2) 26-50% cover'),
    @(256, 43328.52598266204, 'Quadrat survey', 'Florida', 43326.0, 'Aislyn', 14.0, 18.0, 1.9, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(257, 43328.52624133102, 'Quadrat survey', 'Florida', 43326.0, 'Aislyn', 13.0, 0.0, 1.1, 'Silt, Sand, Plant cover', '1) 1-25% cover'),
    @(258, 43328.52650092593, 'Quadrat survey', 'Florida', 43326.0, 'Aislyn', 13.0, 20.0, 1.7, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(259, 43328.5268403588, 'Quadrat survey', 'Florida', 43326.0, 'Aislyn', 12.0, 0.0, 1.1, 'Silt, Sand, Plant cover', '3) 51-75% cover'),
    @(260, 43328.52704613426, 'Quadrat survey', 'Florida', 43326.0, 'Aislyn', 12.0, 5.0, 1.2, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(261, 43328.52724039352, 'Quadrat survey', 'Florida', 43326.0, 'Aislyn', 11.0, 0.0, 1.3, 'Silt, Sand, Plant cover', '4) 75-100% cover'),
    @(262, 43328.527494375, 'Quadrat survey', 'Florida', 43327.0, 'Aislyn', 9.0, 0.0, 1.1, 'Silt, Sand, Gravel, Rock, Plant cover', '1) 1-25% cover'),
    @(263, 43328.52770712963, 'Quadrat survey', 'Florida', 43327.0, 'Aislyn', 10.0, 0.0, 1.3, 'Silt, Sand, Gravel, Plant cover', '1) 1-25% cover'),
    @(264, 43328.52795751157, 'Quadrat survey', 'Florida', 43327.0, 'Aislyn', 10.0, 11.0, 1.6, 'Silt, Sand, Plant cover', '2) 26-50% cover'),
    @(265, 43328.528181435184, 'Quadrat survey', 'Florida', 43327.0, 'Aislyn', 10.0, 24.0, 2.0, 'Silt, Sand, Plant cover', '4) 75-100% cover')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
}

Write-Host ("Added " + $rows.Count + " rows (" + $firstNew + "-" + $lastNew + ") for Lake Florida")
